# Update "想去人数" (want-to-go count, column F) values on the
# "展览" and "全部类型" worksheets to reflect newly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 584
$wsExhibit.Range("F4").Value = 424
$wsExhibit.Range("F5").Value = 465
$wsExhibit.Range("F6").Value = 275
$wsExhibit.Range("F7").Value = 2499
$wsExhibit.Range("F8").Value = 429
$wsExhibit.Range("F9").Value = 6623
$wsExhibit.Range("F11").Value = 425

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 584
$wsAll.Range("F4").Value = 424
$wsAll.Range("F5").Value = 465
$wsAll.Range("F6").Value = 275
$wsAll.Range("F9").Value = 2499
$wsAll.Range("F10").Value = 429
$wsAll.Range("F11").Value = 6623
$wsAll.Range("F13").Value = 425
